# Cambios en el manejo de modelo de factura
#
# The "Datos" sheet holds a small data-driven table used by the BDD tests
# for "InscribirFacturas". The invoice model changed: the old
# empresa/referencia(1) columns (with sample values "movistar" /
# "referencia 1") are replaced by a convenio/referencia pair driven by
# numeric codes ("65401" / "8417").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

# Update the value-row cells first (introduces "65401" then "8417" into the
# shared-string table), then the header cell (introduces "convenio" last) so
# the shared-string insertion order matches the new model.
$ws.Range("M2").Value = "65401"
$ws.Range("O2").Value = "8417"
$ws.Range("M1").Value = "convenio"

# Reflect the new view/selection state left by the edit: scrolled so column E
# is the left-most visible column, with the active cell/selection on M21.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("M21").Select()
